$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.942.39'
$ws.Range("E2").Value = '  -4.77%  '
$ws.Range("D3").Value = '2.493.00'
$ws.Range("E3").Value = '  -3.33%  '
$ws.Range("E4").Value = '  -0.04%  '
$fmt = $ws.Range("D5").NumberFormat
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.15'
$ws.Range("D5").NumberFormat = $fmt
$ws.Range("E5").Value = '  -2.82%  '
$fmt = $ws.Range("D6").NumberFormat
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.89'
$ws.Range("D6").NumberFormat = $fmt
$ws.Range("E6").Value = '  -6.98%  '
$fmt = $ws.Range("D7").NumberFormat
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").NumberFormat = $fmt
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  -3.16%  '
$ws.Range("D9").Value = '2.523.89'
$ws.Range("E9").Value = '  -2.35%  '
$fmt = $ws.Range("D10").NumberFormat
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0998'
$ws.Range("D10").NumberFormat = $fmt
$ws.Range("E10").Value = '  -4.11%  '
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("E12").Value = '  -0.31%  '
$fmt = $ws.Range("D13").NumberFormat
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("D13").NumberFormat = $fmt
$ws.Range("E13").Value = '  -3.79%  '
$ws.Range("D14").Value = '2.930.61'
$ws.Range("E14").Value = '  -3.34%  '
$fmt = $ws.Range("D15").NumberFormat
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.78'
$ws.Range("D15").NumberFormat = $fmt
$ws.Range("E15").Value = '  -6.49%  '
$ws.Range("D16").Value = '58.821.83'
$ws.Range("E16").Value = '  -4.87%  '
$ws.Range("E17").Value = '  -3.87%  '
$ws.Range("D18").Value = '2.511.63'
$ws.Range("E18").Value = '  -2.70%  '
$fmt = $ws.Range("D19").NumberFormat
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.27'
$ws.Range("D19").NumberFormat = $fmt
$ws.Range("E19").Value = '  -2.71%  '
$ws.Range("E20").Value = '  -5.79%  '
$fmt = $ws.Range("D21").NumberFormat
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.80'
$ws.Range("D21").NumberFormat = $fmt
$ws.Range("E21").Value = '  -4.64%  '
$fmt = $ws.Range("D22").NumberFormat
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").NumberFormat = $fmt
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("E23").Value = '  -4.53%  '
$fmt = $ws.Range("D24").NumberFormat
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.66'
$ws.Range("D24").NumberFormat = $fmt
$ws.Range("E24").Value = '  -2.98%  '
$fmt = $ws.Range("D25").NumberFormat
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.439'
$ws.Range("D25").NumberFormat = $fmt
$ws.Range("E25").Value = '  -10.59%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '2.613.07'
$ws.Range("E26").Value = '  -3.00%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$fmt = $ws.Range("D27").NumberFormat
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.161'
$ws.Range("D27").NumberFormat = $fmt
$ws.Range("E27").Value = '  -3.57%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$fmt = $ws.Range("D28").NumberFormat
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.991'
$ws.Range("D28").NumberFormat = $fmt
$ws.Range("E28").Value = '  -0.85%  '
$fmt = $ws.Range("D29").NumberFormat
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.75'
$ws.Range("D29").NumberFormat = $fmt
$ws.Range("E29").Value = '  -4.70%  '
$fmt = $ws.Range("D30").NumberFormat
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.78'
$ws.Range("D30").NumberFormat = $fmt
$ws.Range("E30").Value = '  -6.25%  '
$ws.Range("D31").Value = '0.0₃0775'
$ws.Range("E31").Value = '  -6.92%  '
$ws.Range("E32").Value = '  -5.26%  '
$fmt = $ws.Range("D33").NumberFormat
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.21'
$ws.Range("D33").NumberFormat = $fmt
$ws.Range("E33").Value = '  -10.16%  '
$fmt = $ws.Range("D34").NumberFormat
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '159.74'
$ws.Range("D34").NumberFormat = $fmt
$ws.Range("E34").Value = '  -1.95%  '
$fmt = $ws.Range("D35").NumberFormat
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.995'
$ws.Range("D35").NumberFormat = $fmt
$ws.Range("E35").Value = '  -0.33%  '
$ws.Range("E36").Value = '  +3.32%  '
$fmt = $ws.Range("D37").NumberFormat
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.55'
$ws.Range("D37").NumberFormat = $fmt
$ws.Range("E37").Value = '  -3.15%  '
$fmt = $ws.Range("D38").NumberFormat
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.39'
$ws.Range("D38").NumberFormat = $fmt
$ws.Range("E38").Value = '  -10.08%  '
$fmt = $ws.Range("D39").NumberFormat
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.61'
$ws.Range("D39").NumberFormat = $fmt
$ws.Range("E39").Value = '  -9.50%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$fmt = $ws.Range("D40").NumberFormat
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '306.19'
$ws.Range("D40").NumberFormat = $fmt
$ws.Range("E40").Value = '  -5.94%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$fmt = $ws.Range("D41").NumberFormat
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.61'
$ws.Range("D41").NumberFormat = $fmt
$ws.Range("E41").Value = '  -7.19%  '
$fmt = $ws.Range("D42").NumberFormat
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.78'
$ws.Range("D42").NumberFormat = $fmt
$ws.Range("E42").Value = '  -1.94%  '
$fmt = $ws.Range("D43").NumberFormat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.821'
$ws.Range("D43").NumberFormat = $fmt
$ws.Range("E43").Value = '  -8.78%  '
$fmt = $ws.Range("D44").NumberFormat
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.66'
$ws.Range("D44").NumberFormat = $fmt
$ws.Range("E44").Value = '  -6.55%  '
$fmt = $ws.Range("D45").NumberFormat
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.995'
$ws.Range("D45").NumberFormat = $fmt
$ws.Range("E45").Value = '  -0.34%  '
$fmt = $ws.Range("D46").NumberFormat
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.595'
$ws.Range("D46").NumberFormat = $fmt
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("E47").Value = '  -1.59%  '
$fmt = $ws.Range("D48").NumberFormat
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.71'
$ws.Range("D48").NumberFormat = $fmt
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("E49").Value = '  -3.60%  '
$fmt = $ws.Range("D50").NumberFormat
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.62'
$ws.Range("D50").NumberFormat = $fmt
$ws.Range("E50").Value = '  -4.54%  '
$fmt = $ws.Range("D51").NumberFormat
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0516'
$ws.Range("D51").NumberFormat = $fmt
$ws.Range("E51").Value = '  -5.35%  '
